$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pareto exponent for extrapolation values were incorrectly rounded (left blank).
# Fill in the correct numeric values for columns I (M1_PH), J (CM2_PH),
# K (CMN3_PH), L (CMN4_PH) for rows 12, 13, 14, 16, 18.

$values = @{
    12 = @(-0.1138270110700734, 0.01678756918276458, 0.053487084533311, 2.298063698297935)
    13 = @(-0.08943813976169483, 0.01696372620412041, -0.2880759216279744, 2.14521553204206)
    14 = @(-0.02189067404115775, 0.04191629749799171, 0.6704757719610467, 2.815473210689341)
    16 = @(0.03408971441573483, 0.05861338755323488, 0.4399583757458075, 1.992906231796149)
    18 = @(0.09378818548282215, 0.08061204851819169, 0.3553474179588729, 1.795346708229094)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("I$row").Value = $rowValues[0]
    $ws.Range("J$row").Value = $rowValues[1]
    $ws.Range("K$row").Value = $rowValues[2]
    $ws.Range("L$row").Value = $rowValues[3]
}
